$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.692.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.445.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.444.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("E11").Value = "  -8.56%  "
$ws.Range("E12").Value = "  -7.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.030.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("E14").Value = "  -9.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.444.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.619.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -10.38%  "
$ws.Range("E20").Value = "  -6.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.43%  "
$ws.Range("E23").Value = "  -8.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.587.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("E27").Value = "  -9.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -8.76%  "
$ws.Range("E30").Value = "  -8.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.452.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.147"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.83%  "
$ws.Range("E35").Value = "  -6.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "171.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -8.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.826"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E45").Value = "  -13.84%  "
$ws.Range("E46").Value = "  -11.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.91%  "
$ws.Range("E50").Value = "  -14.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.205.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.54%  "
